# Auto-generated script to apply refreshed market-price data to the Tonberry Profits workbook.
# For each affected sheet, update the computed price/profit columns (H:N) to match the
# latest values pulled by the scheduled market-data runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 119.22222
$ws.Range("I33").Value = 135.42857
$ws.Range("J33").Value = 108.90909
$ws.Range("K33").Value = 135.42857
$ws.Range("L33").Value = 108.90909
$ws.Range("M33").Value = 93.57142999999999
$ws.Range("N33").Value = -566.90909
$ws.Range("H69").Value = 3000
$ws.Range("I69").Value = 3000
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 9000
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -8126
$ws.Range("H72").Value = 3000
$ws.Range("I72").Value = 3000
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 27000
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -22632
$ws.Range("H92").Value = 821137
$ws.Range("I92").Value = 1026254.8
$ws.Range("K92").Value = 1026254.8
$ws.Range("M92").Value = -1025006.8
$ws.Range("H96").Value = 1434.8
$ws.Range("I96").Value = 1434.8
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 4304.4
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -2931.4
$ws.Range("H132").Value = 827.65955
$ws.Range("I132").Value = 766.13635
$ws.Range("J132").Value = 1730
$ws.Range("K132").Value = 2298.40905
$ws.Range("L132").Value = 5190
$ws.Range("M132").Value = 231.5909499999998
$ws.Range("N132").Value = -10250
$ws.Range("H138").Value = 2401.22
$ws.Range("J138").Value = 2342.3333
$ws.Range("L138").Value = 7026.999899999999
$ws.Range("N138").Value = -17306.9999
$ws.Range("N69").ClearContents()
$ws.Range("N72").ClearContents()
$ws.Range("N96").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4888.3
$ws.Range("I61").Value = 1997.7142
$ws.Range("K61").Value = 1997.7142
$ws.Range("M61").Value = -1785.7142
$ws.Range("H74").Value = 1562.421
$ws.Range("I74").Value = 967.2727
$ws.Range("K74").Value = 967.2727
$ws.Range("M74").Value = -93.27269999999999
$ws.Range("H77").Value = 1562.421
$ws.Range("I77").Value = 967.2727
$ws.Range("K77").Value = 4836.363499999999
$ws.Range("M77").Value = -468.3634999999995
$ws.Range("H102").Value = 1799.5
$ws.Range("I102").Value = 1799.5
$ws.Range("K102").Value = 1799.5
$ws.Range("M102").Value = -177.5
$ws.Range("H110").Value = 1688.0416
$ws.Range("I110").Value = 1312.381
$ws.Range("J110").Value = 4317.6665
$ws.Range("K110").Value = 1312.381
$ws.Range("L110").Value = 4317.6665
$ws.Range("M110").Value = 732.6189999999999
$ws.Range("N110").Value = -8407.666499999999
$ws.Range("H136").Value = 4888.3
$ws.Range("I136").Value = 1997.7142
$ws.Range("K136").Value = 5993.142599999999
$ws.Range("M136").Value = -3443.142599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1932.3158
$ws.Range("I20").Value = 1815.7693
$ws.Range("K20").Value = 1815.7693
$ws.Range("M20").Value = -1568.7693
$ws.Range("H81").Value = 18075
$ws.Range("J81").Value = 18075
$ws.Range("L81").Value = 18075
$ws.Range("N81").Value = -20197
$ws.Range("H84").Value = 18075
$ws.Range("J84").Value = 18075
$ws.Range("L84").Value = 54225
$ws.Range("N84").Value = -64833
$ws.Range("H99").Value = 1790.2632
$ws.Range("I99").Value = 1578.0769
$ws.Range("K99").Value = 1578.0769
$ws.Range("M99").Value = -80.07690000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1813589.5
$ws.Range("I58").Value = 2416958.2
$ws.Range("J58").Value = 3482.8333
$ws.Range("K58").Value = 2416958.2
$ws.Range("L58").Value = 3482.8333
$ws.Range("M58").Value = -2416755.2
$ws.Range("N58").Value = -3888.8333
$ws.Range("H132").Value = 2057.9062
$ws.Range("I132").Value = 1172.7142
$ws.Range("J132").Value = 3747.818
$ws.Range("K132").Value = 3518.1426
$ws.Range("L132").Value = 11243.454
$ws.Range("M132").Value = -988.1425999999997
$ws.Range("N132").Value = -16303.454
$ws.Range("H134").Value = 1429.0667
$ws.Range("I134").Value = 1440.4482
$ws.Range("J134").Value = 1099
$ws.Range("K134").Value = 4321.3446
$ws.Range("L134").Value = 3297
$ws.Range("M134").Value = -1786.3446
$ws.Range("N134").Value = -8367
$ws.Range("H136").Value = 1813589.5
$ws.Range("I136").Value = 2416958.2
$ws.Range("J136").Value = 3482.8333
$ws.Range("K136").Value = 7250874.600000001
$ws.Range("L136").Value = 10448.4999
$ws.Range("M136").Value = -7248324.600000001
$ws.Range("N136").Value = -15548.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 404.81818
$ws.Range("I7").Value = 170.2
$ws.Range("J7").Value = 600.3333
$ws.Range("K7").Value = 510.6
$ws.Range("L7").Value = 1800.9999
$ws.Range("M7").Value = -398.6
$ws.Range("N7").Value = -2024.9999
$ws.Range("H92").Value = 318.3
$ws.Range("J92").Value = 332.16666
$ws.Range("L92").Value = 996.4999799999999
$ws.Range("N92").Value = -3492.49998
$ws.Range("H93").Value = 6000
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("H131").Value = 8979.397999999999
$ws.Range("I131").Value = 690
$ws.Range("J131").Value = 9332.138999999999
$ws.Range("K131").Value = 2070
$ws.Range("L131").Value = 27996.417
$ws.Range("M131").Value = 2970
$ws.Range("N131").Value = -38076.417
$ws.Range("H137").Value = 3608.0908
$ws.Range("I137").Value = 2747.8
$ws.Range("J137").Value = 4325
$ws.Range("K137").Value = 8243.400000000001
$ws.Range("L137").Value = 12975
$ws.Range("M137").Value = -3143.400000000001
$ws.Range("N137").Value = -23175
$ws.Range("M93").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 10285
$ws.Range("J109").Value = 10285
$ws.Range("L109").Value = 10285
$ws.Range("N109").Value = -12365

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2482
$ws.Range("I22").Value = 3129.8
$ws.Range("K22").Value = 3129.8
$ws.Range("M22").Value = -2834.8
$ws.Range("H27").Value = 2482
$ws.Range("I27").Value = 3129.8
$ws.Range("K27").Value = 3129.8
$ws.Range("M27").Value = -3022.8
$ws.Range("H40").Value = 4414.778
$ws.Range("I40").Value = 2389.8333
$ws.Range("K40").Value = 2389.8333
$ws.Range("M40").Value = -2253.8333
$ws.Range("H136").Value = 4310.091
$ws.Range("I136").Value = 2882.4443
$ws.Range("J136").Value = 5298.4614
$ws.Range("K136").Value = 8647.332900000001
$ws.Range("L136").Value = 15895.3842
$ws.Range("M136").Value = -6097.332900000001
$ws.Range("N136").Value = -20995.3842

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 497.72726
$ws.Range("I100").Value = 397.5
$ws.Range("J100").Value = 1500
$ws.Range("K100").Value = 795
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -254
$ws.Range("N100").Value = -4082
$ws.Range("H132").Value = 2399.2856
$ws.Range("I132").Value = 1359.4
$ws.Range("K132").Value = 4078.2
$ws.Range("M132").Value = -1548.2

